$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Myoc"
$ws.Range("C2").Value = "Fzd10"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2439023333333333
$ws.Range("H2").Value = 0.731707
$ws.Range("I2").Value = 0.006368708152767561
$ws.Range("J2").Value = 0.006368708152767561
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01234033333333333
$ws.Range("N2").Value = 0.037021
$ws.Range("O2").Value = 0.1810521476743106
$ws.Range("P2").Value = 0.1810521476743105
$ws.Range("Q2").Value = 0.003009836094111111
$ws.Range("R2").Value = 0.027088524847
$ws.Range("S2").Value = 0.001153068288969458
$ws.Range("T2").Value = 0.001153068288969458

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Myoc"
$ws.Range("C3").Value = "Fzd10"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2439023333333333
$ws.Range("H3").Value = 0.731707
$ws.Range("I3").Value = 0.006368708152767561
$ws.Range("J3").Value = 0.006368708152767561
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.05581866666666666
$ws.Range("N3").Value = 0.167456
$ws.Range("O3").Value = 0.8189478523256895
$ws.Range("P3").Value = 0.8189478523256895
$ws.Range("Q3").Value = 0.01361430304355555
$ws.Range("R3").Value = 0.122528727392
$ws.Range("S3").Value = 0.005215639863798104
$ws.Range("T3").Value = 0.005215639863798104

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Myoc"
$ws.Range("C4").Value = "Fzd10"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 37.428193
$ws.Range("H4").Value = 112.284579
$ws.Range("I4").Value = 0.9773142989029397
$ws.Range("J4").Value = 0.9773142989029399
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01234033333333333
$ws.Range("N4").Value = 0.037021
$ws.Range("O4").Value = 0.1810521476743106
$ws.Range("P4").Value = 0.1810521476743105
$ws.Range("Q4").Value = 0.4618763776843333
$ws.Range("R4").Value = 4.156887399159
$ws.Range("S4").Value = 0.1769448527691904
$ws.Range("T4").Value = 0.1769448527691903

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Myoc"
$ws.Range("C5").Value = "Fzd10"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 37.428193
$ws.Range("H5").Value = 112.284579
$ws.Range("I5").Value = 0.9773142989029397
$ws.Range("J5").Value = 0.9773142989029399
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05581866666666666
$ws.Range("N5").Value = 0.167456
$ws.Range("O5").Value = 0.8189478523256895
$ws.Range("P5").Value = 0.8189478523256895
$ws.Range("Q5").Value = 2.089191829002667
$ws.Range("R5").Value = 18.802726461024
$ws.Range("S5").Value = 0.8003694461337494
$ws.Range("T5").Value = 0.8003694461337495

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Myoc"
$ws.Range("C6").Value = "Fzd10"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6248916666666666
$ws.Range("H6").Value = 1.874675
$ws.Range("I6").Value = 0.01631699294429263
$ws.Range("J6").Value = 0.01631699294429263
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01234033333333333
$ws.Range("N6").Value = 0.037021
$ws.Range("O6").Value = 0.1810521476743106
$ws.Range("P6").Value = 0.1810521476743105
$ws.Range("Q6").Value = 0.007711371463888888
$ws.Range("R6").Value = 0.06940234317499999
$ws.Range("S6").Value = 0.002954226616150753
$ws.Range("T6").Value = 0.002954226616150752

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Myoc"
$ws.Range("C7").Value = "Fzd10"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6248916666666666
$ws.Range("H7").Value = 1.874675
$ws.Range("I7").Value = 0.01631699294429263
$ws.Range("J7").Value = 0.01631699294429263
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.05581866666666666
$ws.Range("N7").Value = 0.167456
$ws.Range("O7").Value = 0.8189478523256895
$ws.Range("P7").Value = 0.8189478523256895
$ws.Range("Q7").Value = 0.03488061964444444
$ws.Range("R7").Value = 0.3139255767999999
$ws.Range("S7").Value = 0.01336276632814188
$ws.Range("T7").Value = 0.01336276632814188
